$wb = $excel.ActiveWorkbook

# Work on the "Repayment Schedule" sheet (Variable Instalments columns added)
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late"/"Outstanding" columns one place to the right.
$ws.Columns("N").Insert()

# This sheet becomes the active sheet/tab, with R8 selected.
$ws.Activate()
$ws.Range("R8").Select()
